# chore(results): Auto-update draw results on excel 2025-10-04T17:34:58Z
# Appends the new Pick 4 draw result row (2025-10-04) to the "Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

# Force the new row to be stored as text (matching every other row in the
# table, which keeps values like "2025-10-04" / "251004" from being
# reinterpreted as a date serial / number).
$ws.Range("A" + $row + ":E" + $row).NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-10-04"
$ws.Range("B" + $row).Value = "Pick 4"
$ws.Range("C" + $row).Value = "251004"
$ws.Range("D" + $row).Value = "7-4-5-2"
$ws.Range("E" + $row).Value = "2025-10-04T21:34:58.229+04:00"
